# Add the "Verify" data for the Manager - Bart rows (11-14) in column G,
# and update the active selection as recorded in the saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G11").Value = 28
$ws.Range("G12").Value = 26
$ws.Range("G13").Value = 45
$ws.Range("G14").Value = 48

[void]$ws.Range("H17").Select()
